$wb = $excel.ActiveWorkbook

# --- Update sigma_010 values (B2:C12) ---
$ws2 = $wb.Worksheets.Item("sigma_010")
$ws2.Cells.Item(2, 2).Value = 27.79192028378455
$ws2.Cells.Item(2, 3).Value = 30.05194753110972
$ws2.Cells.Item(3, 2).Value = 27.74342111878395
$ws2.Cells.Item(3, 3).Value = 30.04432959408366
$ws2.Cells.Item(4, 2).Value = 27.80352253831376
$ws2.Cells.Item(4, 3).Value = 30.04787804451952
$ws2.Cells.Item(5, 2).Value = 27.77638544348756
$ws2.Cells.Item(5, 3).Value = 30.05361456135043
$ws2.Cells.Item(6, 2).Value = 27.81568648223155
$ws2.Cells.Item(6, 3).Value = 30.07292473239581
$ws2.Cells.Item(7, 2).Value = 27.77566706697283
$ws2.Cells.Item(7, 3).Value = 30.04532183104329
$ws2.Cells.Item(8, 2).Value = 27.78064485728103
$ws2.Cells.Item(8, 3).Value = 30.02379336207301
$ws2.Cells.Item(9, 2).Value = 27.76388528385115
$ws2.Cells.Item(9, 3).Value = 30.03525227144778
$ws2.Cells.Item(10, 2).Value = 27.7606223056024
$ws2.Cells.Item(10, 3).Value = 30.05166986494365
$ws2.Cells.Item(11, 2).Value = 27.81377968849252
$ws2.Cells.Item(11, 3).Value = 30.06028049097981
$ws2.Cells.Item(12, 2).Value = 27.78255350688013
$ws2.Cells.Item(12, 3).Value = 30.04870122839467

# --- Update sigma_025 values (B2:C12) ---
$ws3 = $wb.Worksheets.Item("sigma_025")
$ws3.Cells.Item(2, 2).Value = 19.65849188845619
$ws3.Cells.Item(2, 3).Value = 25.84959761299686
$ws3.Cells.Item(3, 2).Value = 19.67991510124479
$ws3.Cells.Item(3, 3).Value = 25.80411902149311
$ws3.Cells.Item(4, 2).Value = 19.66376740704578
$ws3.Cells.Item(4, 3).Value = 25.78439205124278
$ws3.Cells.Item(5, 2).Value = 19.63920777154176
$ws3.Cells.Item(5, 3).Value = 25.82807974025393
$ws3.Cells.Item(6, 2).Value = 19.66933018534947
$ws3.Cells.Item(6, 3).Value = 25.85284403072756
$ws3.Cells.Item(7, 2).Value = 19.65868365857116
$ws3.Cells.Item(7, 3).Value = 25.83254329819898
$ws3.Cells.Item(8, 2).Value = 19.65453909434124
$ws3.Cells.Item(8, 3).Value = 25.81550334835927
$ws3.Cells.Item(9, 2).Value = 19.67448632700927
$ws3.Cells.Item(9, 3).Value = 25.84372877560067
$ws3.Cells.Item(10, 2).Value = 19.68461255960584
$ws3.Cells.Item(10, 3).Value = 25.83938499894241
$ws3.Cells.Item(11, 2).Value = 19.65208493996195
$ws3.Cells.Item(11, 3).Value = 25.81814612833216
$ws3.Cells.Item(12, 2).Value = 19.66351189331274
$ws3.Cells.Item(12, 3).Value = 25.82683390061477

# --- Add new sheet sigma_050 at the end ---
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws4.Name = "sigma_050"

# --- Populate header row ---
$ws4.Cells.Item(1, 1).Value = "Rows"
$ws4.Cells.Item(1, 2).Value = "Noisy"
$ws4.Cells.Item(1, 3).Value = "NLM-LBP"

# --- Populate data rows (A: index, B: Noisy, C: NLM-LBP) ---
$ws4.Cells.Item(2, 1).Value = 0
$ws4.Cells.Item(2, 2).Value = 14.58870232807314
$ws4.Cells.Item(2, 3).Value = 21.51692180155847
$ws4.Cells.Item(3, 1).Value = 1
$ws4.Cells.Item(3, 2).Value = 14.57381268532069
$ws4.Cells.Item(3, 3).Value = 21.53141173778381
$ws4.Cells.Item(4, 1).Value = 2
$ws4.Cells.Item(4, 2).Value = 14.58893345443682
$ws4.Cells.Item(4, 3).Value = 21.47661076888775
$ws4.Cells.Item(5, 1).Value = 3
$ws4.Cells.Item(5, 2).Value = 14.608672509545
$ws4.Cells.Item(5, 3).Value = 21.55451958826844
$ws4.Cells.Item(6, 1).Value = 4
$ws4.Cells.Item(6, 2).Value = 14.60293346479738
$ws4.Cells.Item(6, 3).Value = 21.5834265626049
$ws4.Cells.Item(7, 1).Value = 5
$ws4.Cells.Item(7, 2).Value = 14.59997342616882
$ws4.Cells.Item(7, 3).Value = 21.58188092387499
$ws4.Cells.Item(8, 1).Value = 6
$ws4.Cells.Item(8, 2).Value = 14.61558380071422
$ws4.Cells.Item(8, 3).Value = 21.58421078214695
$ws4.Cells.Item(9, 1).Value = 7
$ws4.Cells.Item(9, 2).Value = 14.59726110953686
$ws4.Cells.Item(9, 3).Value = 21.52919309973935
$ws4.Cells.Item(10, 1).Value = 8
$ws4.Cells.Item(10, 2).Value = 14.5967687029124
$ws4.Cells.Item(10, 3).Value = 21.54595098752557
$ws4.Cells.Item(11, 1).Value = 9
$ws4.Cells.Item(11, 2).Value = 14.59301333679739
$ws4.Cells.Item(11, 3).Value = 21.5252725731062

# --- Final "Média" (average) row ---
$ws4.Cells.Item(12, 1).Value = "Média"
$ws4.Cells.Item(12, 2).Value = 14.59656548183027
$ws4.Cells.Item(12, 3).Value = 21.54293988254964
